# Update "想去人数" (column F) values across the four sheets to match
# the freshly generated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 647
$ws.Cells.Item(3, 6).Value = 728
$ws.Cells.Item(4, 6).Value = 950
$ws.Cells.Item(5, 6).Value = 733
$ws.Cells.Item(6, 6).Value = 841
$ws.Cells.Item(8, 6).Value = 615
$ws.Cells.Item(9, 6).Value = 135
$ws.Cells.Item(10, 6).Value = 1221
$ws.Cells.Item(11, 6).Value = 646
$ws.Cells.Item(12, 6).Value = 389
$ws.Cells.Item(13, 6).Value = 515
$ws.Cells.Item(15, 6).Value = 15
$ws.Cells.Item(16, 6).Value = 576
$ws.Cells.Item(18, 6).Value = 361
$ws.Cells.Item(19, 6).Value = 355
$ws.Cells.Item(22, 6).Value = 87
$ws.Cells.Item(23, 6).Value = 585
$ws.Cells.Item(24, 6).Value = 29
$ws.Cells.Item(25, 6).Value = 802

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 88
$ws.Cells.Item(4, 6).Value = 325
$ws.Cells.Item(5, 6).Value = 105
$ws.Cells.Item(9, 6).Value = 225
$ws.Cells.Item(10, 6).Value = 50
$ws.Cells.Item(13, 6).Value = 100

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 369

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 369
$ws.Cells.Item(3, 6).Value = 88
$ws.Cells.Item(4, 6).Value = 647
$ws.Cells.Item(6, 6).Value = 325
$ws.Cells.Item(7, 6).Value = 728
$ws.Cells.Item(8, 6).Value = 950
$ws.Cells.Item(9, 6).Value = 733
$ws.Cells.Item(10, 6).Value = 841
$ws.Cells.Item(12, 6).Value = 615
$ws.Cells.Item(13, 6).Value = 135
$ws.Cells.Item(14, 6).Value = 1221
$ws.Cells.Item(15, 6).Value = 646
$ws.Cells.Item(16, 6).Value = 105
$ws.Cells.Item(18, 6).Value = 389
$ws.Cells.Item(19, 6).Value = 515
$ws.Cells.Item(22, 6).Value = 15
$ws.Cells.Item(23, 6).Value = 576
$ws.Cells.Item(26, 6).Value = 361
$ws.Cells.Item(27, 6).Value = 355
$ws.Cells.Item(29, 6).Value = 225
$ws.Cells.Item(30, 6).Value = 50
$ws.Cells.Item(34, 6).Value = 100
$ws.Cells.Item(35, 6).Value = 100
$ws.Cells.Item(36, 6).Value = 87
$ws.Cells.Item(37, 6).Value = 585
$ws.Cells.Item(38, 6).Value = 29
$ws.Cells.Item(39, 6).Value = 802

$wb.Save()
